$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy columns N,O,AM,AN (rows 1-3) into columns B,C,D,E (rows 1-3)
$ws.Range("B1:B3").Value2 = $ws.Range("N1:N3").Value2
$ws.Range("C1:C3").Value2 = $ws.Range("O1:O3").Value2
$ws.Range("D1:D3").Value2 = $ws.Range("AM1:AM3").Value2
$ws.Range("E1:E3").Value2 = $ws.Range("AN1:AN3").Value2

# Update the selected range to reflect the new selection (B1:E3)
$ws.Range("B1:E3").Select()
